$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Structural changes -----------------------------------------------
# Insert a new column before "O" (Spcl Allowance) to make room for the
# new "Education Allowance" column. This shifts O..X (and beyond) right
# by one.
$ws.Range("O1").EntireColumn.Insert()

# The old "ABC" column (originally V1) has now been shifted to W1; it is
# dropped entirely in the new layout, so delete that column. This shifts
# total_deducations/NetPay back down into W/X.
$ws.Range("W1").EntireColumn.Delete()

# Remove the last data row (pattabhi ramarao) entirely.
$ws.Rows("5:5").Delete()

# --- Header row (row 1) -------------------------------------------------
$ws.Cells.Item(1, 14).Value = "Car Allowance"        # N1 (was "ABC1")
$ws.Cells.Item(1, 15).Value = "Education Allowance"  # O1 (new column)

# --- Row 2 (Sekhar Beri) -------------------------------------------------
$ws.Cells.Item(2, 6).Value  = "Junior Development"   # F2 DESIGNATION
$ws.Cells.Item(2, 7).Value  = "Development"          # G2 DEPARTMENT
$ws.Cells.Item(2, 8).Value  = 18000.0                # H2 GROSS
$ws.Cells.Item(2, 9).Value  = 1500.0                 # I2 Per Month
$ws.Cells.Item(2, 10).Value = 1401.57                # J2 Actual Per Month
$ws.Cells.Item(2, 13).Value = 601.5                  # M2 BASIC
$ws.Cells.Item(2, 14).Value = 25.0                   # N2 Car Allowance
$ws.Cells.Item(2, 15).Value = 60.15                  # O2 Education Allowance
$ws.Cells.Item(2, 16).Value = 714.92                 # P2 Spcl Allowance
$ws.Cells.Item(2, 17).Value = 0.0                    # Q2 Arrears
$ws.Cells.Item(2, 18).Value = 1401.57                # R2 Gross Pay
$ws.Cells.Item(2, 19).Value = 72.18                  # S2 PF
$ws.Cells.Item(2, 20).Value = 24.53                  # T2 ESIC
$ws.Cells.Item(2, 21).Value = 0.0                    # U2 PT
$ws.Cells.Item(2, 22).Value = 0.0                    # V2 TDS
$ws.Cells.Item(2, 23).Value = 96.71                  # W2 total_deducations
$ws.Cells.Item(2, 24).Value = 1304.86                # X2 NetPay

# --- Row 3 (BalaRaju Vankala) --------------------------------------------
$ws.Cells.Item(3, 8).Value  = 20000.0                 # H3 GROSS
$ws.Cells.Item(3, 9).Value  = 1666.6666666666667      # I3 Per Month
$ws.Cells.Item(3, 10).Value = 1557.3                  # J3 Actual Per Month
$ws.Cells.Item(3, 13).Value = 668.333                 # M3 BASIC
$ws.Cells.Item(3, 14).Value = 25.0                    # N3 Car Allowance
$ws.Cells.Item(3, 15).Value = 0                       # O3 Education Allowance
$ws.Cells.Item(3, 16).Value = 863.97                  # P3 Spcl Allowance
$ws.Cells.Item(3, 17).Value = 0.0                     # Q3 Arrears
$ws.Cells.Item(3, 18).Value = 1557.3                  # R3 Gross Pay
$ws.Cells.Item(3, 19).Value = 80.2                    # S3 PF
$ws.Cells.Item(3, 20).Value = 27.25                   # T3 ESIC
$ws.Cells.Item(3, 21).Value = 0.0                     # U3 PT
$ws.Cells.Item(3, 22).Value = 0.0                     # V3 TDS
$ws.Cells.Item(3, 23).Value = 107.45                  # W3 total_deducations
$ws.Cells.Item(3, 24).Value = 1449.85                 # X3 NetPay

# --- Row 4 (Priyanka Muddana) --------------------------------------------
$ws.Cells.Item(4, 8).Value  = 235235.0                # H4 GROSS
$ws.Cells.Item(4, 9).Value  = 19602.916666666668      # I4 Per Month
$ws.Cells.Item(4, 10).Value = 18318.9                 # J4 Actual Per Month
$ws.Cells.Item(4, 13).Value = 7841.17                 # M4 BASIC
$ws.Cells.Item(4, 15).Value = 0                       # O4 Education Allowance
$ws.Cells.Item(4, 16).Value = 10477.8                 # P4 Spcl Allowance
$ws.Cells.Item(4, 17).Value = 0.0                     # Q4 Arrears
$ws.Cells.Item(4, 18).Value = 18318.9                 # R4 Gross Pay
$ws.Cells.Item(4, 19).Value = 940.94                  # S4 PF
$ws.Cells.Item(4, 20).Value = 320.58                  # T4 ESIC
$ws.Cells.Item(4, 21).Value = 0.0                     # U4 PT
$ws.Cells.Item(4, 22).Value = 0.0                     # V4 TDS
$ws.Cells.Item(4, 23).Value = 1261.52                 # W4 total_deducations
$ws.Cells.Item(4, 24).Value = 17057.4                 # X4 NetPay

# --- Column widths --------------------------------------------------------
# ColumnWidth is expressed in "characters"; the host quantizes to 1/7-char
# pixel steps, so feed values chosen to land as close as possible to the
# authored widths.
$ws.Columns(6).ColumnWidth  = 16.714285714285715   # F  -> 17.38988764044944
$ws.Columns(9).ColumnWidth  = 21.142857142857142   # I  -> 21.789887640449443
$ws.Columns(13).ColumnWidth = 9.0                  # M  -> 9.68988764044944
$ws.Columns(14).ColumnWidth = 10.142857142857142   # N  -> 10.78988764044944
$ws.Columns(15).ColumnWidth = 14.428571428571429   # O  -> 15.18988764044944
$ws.Columns(16).ColumnWidth = 11.142857142857142   # P  -> 11.88988764044944
$ws.Columns(17).ColumnWidth = 4.571428571428571    # Q  -> 5.289887640449439
$ws.Columns(18).ColumnWidth = 9.0                  # R  -> 9.68988764044944
$ws.Columns(20).ColumnWidth = 7.857142857142857    # T  -> 8.589887640449438
$ws.Columns(21).ColumnWidth = 4.571428571428571    # U  -> 5.289887640449439
